$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 / column B ("R40" in the rules table) is retyped as the text "1".
# The leading apostrophe forces Excel to store it as text (a new shared
# string "1") instead of auto-converting the numeric-looking value to a
# number, matching how this was entered through the UI.
$ws.Range("B11").Value = "`'1"
